$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at B:C (shifts old B.. onward to D..)
$ws.Range("B1:C1").EntireColumn.Insert()

# Populate the two new header cells
$ws.Range("B1").Value = "Sport"
$ws.Range("C1").Value = "Team"

# Restore the active selection to C2, as recorded in the saved workbook
$ws.Range("C2").Select() | Out-Null
